# Update the "want to go" counts (column F) on the "展览" and "全部类型"
# sheets to reflect refreshed data generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 65
$ws1.Range("F3").Value  = 11658
$ws1.Range("F5").Value  = 336
$ws1.Range("F7").Value  = 11622
$ws1.Range("F9").Value  = 1165
$ws1.Range("F10").Value = 91
$ws1.Range("F12").Value = 5749
$ws1.Range("F13").Value = 117
$ws1.Range("F14").Value = 3511

# --- Sheet "全部类型" ---------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 65
$ws4.Range("F5").Value  = 11658
$ws4.Range("F7").Value  = 336
$ws4.Range("F9").Value  = 11622
$ws4.Range("F11").Value = 1165
$ws4.Range("F12").Value = 92
$ws4.Range("F15").Value = 5749
$ws4.Range("F16").Value = 117
$ws4.Range("F17").Value = 3511
